$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2936666666666667
$ws.Range("N2").Value = 0.881
$ws.Range("O2").Value = 0.009113820319201367
$ws.Range("P2").Value = 0.009113820319201367
$ws.Range("Q2").Value = 13.869982974
$ws.Range("R2").Value = 124.829846766
$ws.Range("S2").Value = 0.002956752406031069
$ws.Range("T2").Value = 0.002956752406031069
# Row 3
$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.870405726797791
$ws.Range("P3").Value = 0.870405726797791
$ws.Range("Q3").Value = 1324.637988058926
$ws.Range("R3").Value = 11921.74189253033
$ws.Range("S3").Value = 0.2823814971983241
$ws.Range("T3").Value = 0.2823814971983241
# Row 4
$ws.Range("G4").Value = 47.23036199999999
$ws.Range("H4").Value = 141.691086
$ws.Range("I4").Value = 0.3244251370417807
$ws.Range("J4").Value = 0.3244251370417807
$ws.Range("O4").Value = 0.1204804528830076
$ws.Range("P4").Value = 0.1204804528830076
$ws.Range("Q4").Value = 183.354704356686
$ws.Range("R4").Value = 1650.192339210174
$ws.Range("S4").Value = 0.03908688743742553
$ws.Range("T4").Value = 0.03908688743742553
# Row 5
$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2936666666666667
$ws.Range("N5").Value = 0.881
$ws.Range("O5").Value = 0.009113820319201367
$ws.Range("P5").Value = 0.009113820319201367
$ws.Range("Q5").Value = 17.90709039322222
$ws.Range("R5").Value = 161.163813539
$ws.Range("S5").Value = 0.003817368247994769
$ws.Range("T5").Value = 0.003817368247994769
# Row 6
$ws.Range("I6").Value = 0.4188548944674916
$ws.Range("J6").Value = 0.4188548944674916
$ws.Range("O6").Value = 0.870405726797791
$ws.Range("P6").Value = 0.870405726797791
$ws.Range("S6").Value = 0.3645736988417891
$ws.Range("T6").Value = 0.3645736988417891
# Row 7
$ws.Range("I7").Value = 0.4188548944674916
$ws.Range("J7").Value = 0.4188548944674916
$ws.Range("O7").Value = 0.1204804528830076
$ws.Range("P7").Value = 0.1204804528830076
$ws.Range("S7").Value = 0.05046382737770773
$ws.Range("T7").Value = 0.05046382737770773
# Row 8
$ws.Range("I8").Value = 0.2567199684907278
$ws.Range("J8").Value = 0.2567199684907277
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2936666666666667
$ws.Range("N8").Value = 0.881
$ws.Range("O8").Value = 0.009113820319201367
$ws.Range("P8").Value = 0.009113820319201367
$ws.Range("Q8").Value = 10.97541831844445
$ws.Range("R8").Value = 98.778764866
$ws.Range("S8").Value = 0.002339699665175529
$ws.Range("T8").Value = 0.002339699665175529
# Row 9
$ws.Range("I9").Value = 0.2567199684907278
$ws.Range("J9").Value = 0.2567199684907277
$ws.Range("O9").Value = 0.870405726797791
$ws.Range("P9").Value = 0.870405726797791
$ws.Range("S9").Value = 0.2234505307576779
$ws.Range("T9").Value = 0.2234505307576779
# Row 10
$ws.Range("I10").Value = 0.2567199684907278
$ws.Range("J10").Value = 0.2567199684907277
$ws.Range("O10").Value = 0.1204804528830076
$ws.Range("P10").Value = 0.1204804528830076
$ws.Range("S10").Value = 0.03092973806787432
$ws.Range("T10").Value = 0.03092973806787431
